$d = $word.ActiveDocument
$wmain = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Hunk 1: "CU-03 Registrar producto" -> add a second run
#         " extendido de CU-04 Ver producto"
# ---------------------------------------------------------------------------
$r1 = $d.Content
$ok1 = $r1.Find.Execute("CU-03 Registrar producto", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok1) {
    $xml1 = '<w:p ' + $wmain + '><w:pPr><w:jc w:val="both"/></w:pPr>' + `
        '<w:r><w:t>CU-03 Registrar producto</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> extendido de CU-04 Ver producto</w:t></w:r>' + `
        '</w:p>'
    [void]$r1.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# Hunk 2: "El actor da clic en el boton "Registrar". (FA-01, FA-02)"
#          -> "... (FA-01,FA-02)" with proofErr gramStart/gramEnd around "01,FA"
# ---------------------------------------------------------------------------
$r2 = $d.Content
$needle2 = "El actor da clic en el bot" + [char]0x00F3 + "n " + [char]0x201C + "Registrar" + [char]0x201D + ". (FA-01, FA-02)"
$ok2 = $r2.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok2) {
    $txt2 = "El actor da clic en el bot" + [char]0x00F3 + "n " + [char]0x201C + "Registrar" + [char]0x201D + ". (FA-"
    $xml2 = '<w:p ' + $wmain + '><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:jc w:val="both"/></w:pPr>' + `
        '<w:r><w:t>' + $txt2 + '</w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>01,FA</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t>-02)</w:t></w:r>' + `
        '</w:p>'
    [void]$r2.InsertXML($xml2)
}

# ---------------------------------------------------------------------------
# Hunk 3: split sentence to insert "cada productos" wrapped with
#         proofErr gramStart/gramEnd
# ---------------------------------------------------------------------------
$r3 = $d.Content
$needle3 = "El sistema verifica que los datos registrados en cada producto sean correctos, luego, el sistema modifica los PRODUCTOs del PEDIDO con el estado de ubicaci" + [char]0x00F3 + "n " + [char]0x201C + "En inventario" + [char]0x201D + " asocia la CATEGORIA seleccionada de cada PRODUCTO, tambi" + [char]0x00E9 + "n agrega el estado " + [char]0x201C + "Entregado" + [char]0x201D + " y la fecha de entrega del PEDIDO seleccionado, tambi" + [char]0x00E9 + "n se registra el precio de venta actual de cada PRODUCTO sumando la ganancia registrada con el precio de compra del PRODUCTO en la base de datos (EX-01),"
$ok3 = $r3.Find.Execute($needle3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok3) {
    $tail3 = " del PEDIDO con el estado de ubicaci" + [char]0x00F3 + "n " + [char]0x201C + "En inventario" + [char]0x201D + " asocia la CATEGORIA seleccionada de cada PRODUCTO, tambi" + [char]0x00E9 + "n agrega el estado " + [char]0x201C + "Entregado" + [char]0x201D + " y la fecha de entrega del PEDIDO seleccionado, tambi" + [char]0x00E9 + "n se registra el precio de venta actual de cada PRODUCTO sumando la ganancia registrada con el precio de compra del PRODUCTO en la base de datos (EX-01),"
    $xml3 = '<w:p ' + $wmain + '><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:jc w:val="both"/></w:pPr>' + `
        '<w:r><w:t xml:space="preserve">El sistema verifica que los datos registrados en </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>cada productos</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> sean correctos, luego, el sistema modifica los </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:t>PRODUCTOs</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve">' + $tail3 + '</w:t></w:r>' + `
        '</w:p>'
    [void]$r3.InsertXML($xml3)
}

# ---------------------------------------------------------------------------
# Hunk 4: new table row "Extiende" / "CU-04" appended after "Postcondiciones"
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Extiende"
$newRow.Cells.Item(2).Range.Text = "CU-04"

Write-Host "edits applied"
